$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46042, 5607, 4256, 3885, 235, 72, 52, 9, 3),
    @(46044, 5601, 4399, 3879, 305, 106, 87, 20, 2),
    @(46043, 5602, 4389, 3067, 435, 370, 419, 90, 8)
)

$startRow = 73
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(72, 1).NumberFormat
}

$ws.Range("A75:I75").Select()
